# Burndown chart updated as result of the daily scrum meeting today.
#
# Sprint_One: Philip Guy logged 1 story point done on day M (2017-02-17).
# Enter the value in M2, which mirrors exactly what happens in Excel when a
# user selects M2, types 1 and hits Enter (the selection lands on M3
# afterwards). All of the dependent SUM/running-total formulas in row 6, 7
# (and the Velocity_Chart summary that subtracts Sprint_One!W7 from B7)
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint_One")

# Make sure we are working on the right sheet (it was already the active one).
$ws.Activate()

# Select the cell the user is about to edit, just like a live session would.
$ws.Range("M2").Select()
$ws.Range("M2").Value = 1

# After typing a value and pressing Enter, Excel moves the active cell down
# one row - reproduce that resulting selection state.
$ws.Range("M3").Select()
